# Update the "Latest period (release date)" footnote for the
# "Job adverts by occupation" row on the data table sheet, and move the
# active selection from B12 to A6 (as left by the author after editing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data sources table / job ads footnote: Nov 2024 (07/02/25) -> Jan 2025 (12/03/25)
$ws.Range("C13").Value = "Jan 2025 (12/03/25)"

# Move the selection as recorded in the saved workbook view.
$ws.Range("A6").Select()
